$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51 (Coin, Link, Price, Volume(1h))
$data = @(
    ,@('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '66.144.12', '  -0.59%  ')
    ,@('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.317.86', '  -0.38%  ')
    ,@('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.02%  ')
    ,@('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '586.18', '  +2.34%  ')
    ,@('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '183.26', '  +1.36%  ')
    ,@('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.650', '  +4.72%  ')
    ,@('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.03%  ')
    ,@('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.315.04', '  -0.45%  ')
    ,@('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.126', '  -2.95%  ')
    ,@('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '6.81', '  +2.40%  ')
    ,@('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.401', '  +0.01%  ')
    ,@('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.893.14', '  -0.42%  ')
    ,@('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.132', '  -2.83%  ')
    ,@('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '66.180.96', '  -0.68%  ')
    ,@('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '26.13', '  -2.79%  ')
    ,@('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000163', '  -2.43%  ')
    ,@('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.246.92', '  -2.65%  ')
    ,@('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '424.84', '  -2.69%  ')
    ,@('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.54', '  -2.25%  ')
    ,@('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '13.16', '  -2.60%  ')
    ,@('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '7.39', '  -2.52%  ')
    ,@('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '71.71', '  -2.45%  ')
    ,@('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.02%  ')
    ,@('LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '5.68', '  +0.43%  ')
    ,@('WrappedeETH', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth', '3.463.72', '  -0.80%  ')
    ,@('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.513', '  -0.75%  ')
    ,@('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.200', '  +4.91%  ')
    ,@('PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0000113', '  -3.24%  ')
    ,@('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '8.93', '  -0.95%  ')
    ,@('Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '1.00', '  +0.34%  ')
    ,@('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.91', '  -2.56%  ')
    ,@('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '22.41', '  -1.62%  ')
    ,@('USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  +0.09%  ')
    ,@('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '5.17', '  -2.28%  ')
    ,@('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.56', '  -3.21%  ')
    ,@('Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.18', '  -3.77%  ')
    ,@('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '160.78', '  -1.17%  ')
    ,@('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.44', '  -3.09%  ')
    ,@('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.882.19', '  +2.58%  ')
    ,@('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.80', '  -1.76%  ')
    ,@('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '26.39', '  -3.49%  ')
    ,@('Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.763', '  -3.49%  ')
    ,@('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.30', '  -2.54%  ')
    ,@('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '39.93', '  -0.54%  ')
    ,@('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0663', '  -0.32%  ')
    ,@('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '5.91', '  -4.84%  ')
    ,@('dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.30', '  -2.58%  ')
    ,@('InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '23.14', '  -5.00%  ')
    ,@('Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '313.59', '  -1.82%  ')
    ,@('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0271', '  -0.43%  ')
)

$lastRow = 1 + $data.Count
$ws.Range("D2:D$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}

$wb.Save()